$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2457683333333333
$ws.Range("H2").Value = 0.737305
$ws.Range("I2").Value = 0.1447271191911903
$ws.Range("J2").Value = 0.1575855905380038
$ws.Range("M2").Value = 12.997753
$ws.Range("N2").Value = 38.993259
$ws.Range("O2").Value = 0.4740421406233454
$ws.Range("P2").Value = 0.5546503645614554
$ws.Range("Q2").Value = 3.194436091888334
$ws.Range("R2").Value = 28.749924826995
$ws.Range("S2").Value = 0.06860675338764192
$ws.Range("T2").Value = 0.08740490524153602

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2457683333333333
$ws.Range("H3").Value = 0.737305
$ws.Range("I3").Value = 0.1447271191911903
$ws.Range("J3").Value = 0.1575855905380038
$ws.Range("O3").Value = 0.02725306609819269
$ws.Range("P3").Value = 0.03188729809316786
$ws.Range("Q3").Value = 0.1836507147744444
$ws.Range("R3").Value = 1.65285643297
$ws.Range("S3").Value = 0.003944257745518522
$ws.Range("T3").Value = 0.005024978700673218

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2457683333333333
$ws.Range("H4").Value = 0.737305
$ws.Range("I4").Value = 0.1447271191911903
$ws.Range("J4").Value = 0.1575855905380038
$ws.Range("M4").Value = 1.182591666666666
$ws.Range("N4").Value = 3.547775
$ws.Range("O4").Value = 0.04313039993528083
$ws.Range("P4").Value = 0.05046448405689858
$ws.Range("Q4").Value = 0.2906435829305555
$ws.Range("R4").Value = 2.615792246375
$ws.Range("S4").Value = 0.006242138532197097
$ws.Range("T4").Value = 0.007952475521302039

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2457683333333333
$ws.Range("H5").Value = 0.737305
$ws.Range("I5").Value = 0.1447271191911903
$ws.Range("J5").Value = 0.1575855905380038
$ws.Range("M5").Value = 11.9545335
$ws.Range("N5").Value = 23.909067
$ws.Range("O5").Value = 0.4359947946767024
$ws.Range("P5").Value = 0.3400888529957002
$ws.Range("Q5").Value = 2.9380457740725
$ws.Range("R5").Value = 17.628274644435
$ws.Range("S5").Value = 0.06310027061591367
$ws.Range("T5").Value = 0.05359310273471976

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2457683333333333
$ws.Range("H6").Value = 0.737305
$ws.Range("I6").Value = 0.1447271191911903
$ws.Range("J6").Value = 0.1575855905380038
$ws.Range("M6").Value = 0.5368526666666666
$ws.Range("N6").Value = 1.610558
$ws.Range("O6").Value = 0.01957959866647858
$ws.Range("P6").Value = 0.022909000292778
$ws.Range("Q6").Value = 0.1319413851322222
$ws.Range("R6").Value = 1.18747246619
$ws.Range("S6").Value = 0.002833698909919116
$ws.Range("T6").Value = 0.003610128339772722

# Row 7
$ws.Range("I7").Value = 0.6104826439049008
$ws.Range("J7").Value = 0.6647217777192627
$ws.Range("M7").Value = 12.997753
$ws.Range("N7").Value = 38.993259
$ws.Range("O7").Value = 0.4740421406233454
$ws.Range("P7").Value = 0.5546503645614554
$ws.Range("Q7").Value = 13.474653555323
$ws.Range("R7").Value = 121.271881997907
$ws.Range("S7").Value = 0.2893944993300787
$ws.Range("T7").Value = 0.3686881763439278

# Row 8
$ws.Range("I8").Value = 0.6104826439049008
$ws.Range("J8").Value = 0.6647217777192627
$ws.Range("O8").Value = 0.02725306609819269
$ws.Range("P8").Value = 0.03188729809316786
$ws.Range("S8").Value = 0.01663752384613969
$ws.Range("T8").Value = 0.0211961814751546

# Row 9
$ws.Range("I9").Value = 0.6104826439049008
$ws.Range("J9").Value = 0.6647217777192627
$ws.Range("M9").Value = 1.182591666666666
$ws.Range("N9").Value = 3.547775
$ws.Range("O9").Value = 0.04313039993528083
$ws.Range("P9").Value = 0.05046448405689858
$ws.Range("Q9").Value = 1.225982137508333
$ws.Range("R9").Value = 11.033839237575
$ws.Range("S9").Value = 0.02633036058516601
$ws.Range("T9").Value = 0.03354484155398701

# Row 10
$ws.Range("I10").Value = 0.6104826439049008
$ws.Range("J10").Value = 0.6647217777192627
$ws.Range("M10").Value = 11.9545335
$ws.Range("N10").Value = 23.909067
$ws.Range("O10").Value = 0.4359947946767024
$ws.Range("P10").Value = 0.3400888529957002
$ws.Range("Q10").Value = 12.3931572886485
$ws.Range("R10").Value = 74.358943731891
$ws.Range("S10").Value = 0.2661672549830077
$ws.Range("T10").Value = 0.2260644669458068

# Row 11
$ws.Range("I11").Value = 0.6104826439049008
$ws.Range("J11").Value = 0.6647217777192627
$ws.Range("M11").Value = 0.5368526666666666
$ws.Range("N11").Value = 1.610558
$ws.Range("O11").Value = 0.01957959866647858
$ws.Range("P11").Value = 0.022909000292778
$ws.Range("Q11").Value = 0.5565503278593333
$ws.Range("R11").Value = 5.008952950734
$ws.Range("S11").Value = 0.01195300516050871
$ws.Range("T11").Value = 0.0152281114003865

# Row 12
$ws.Range("G12").Value = 0.4156905
$ws.Range("H12").Value = 0.831381
$ws.Range("I12").Value = 0.2447902369039089
$ws.Range("J12").Value = 0.1776926317427335
$ws.Range("M12").Value = 12.997753
$ws.Range("N12").Value = 38.993259
$ws.Range("O12").Value = 0.4740421406233454
$ws.Range("P12").Value = 0.5546503645614554
$ws.Range("Q12").Value = 5.403042443446501
$ws.Range("R12").Value = 32.418254660679
$ws.Range("S12").Value = 0.1160408879056248
$ws.Range("T12").Value = 0.09855728297599158

# Row 13
$ws.Range("G13").Value = 0.4156905
$ws.Range("H13").Value = 0.831381
$ws.Range("I13").Value = 0.2447902369039089
$ws.Range("J13").Value = 0.1776926317427335
$ws.Range("O13").Value = 0.02725306609819269
$ws.Range("P13").Value = 0.03188729809316786
$ws.Range("Q13").Value = 0.310625280379
$ws.Range("R13").Value = 1.863751682274
$ws.Range("S13").Value = 0.006671284506534475
$ws.Range("T13").Value = 0.005666137917340045

# Row 14
$ws.Range("G14").Value = 0.4156905
$ws.Range("H14").Value = 0.831381
$ws.Range("I14").Value = 0.2447902369039089
$ws.Range("J14").Value = 0.1776926317427335
$ws.Range("M14").Value = 1.182591666666666
$ws.Range("N14").Value = 3.547775
$ws.Range("O14").Value = 0.04313039993528083
$ws.Range("P14").Value = 0.05046448405689858
$ws.Range("Q14").Value = 0.4915921212125
$ws.Range("R14").Value = 2.949552727275
$ws.Range("S14").Value = 0.01055790081791773
$ws.Range("T14").Value = 0.008967166981609523

# Row 15
$ws.Range("G15").Value = 0.4156905
$ws.Range("H15").Value = 0.831381
$ws.Range("I15").Value = 0.2447902369039089
$ws.Range("J15").Value = 0.1776926317427335
$ws.Range("M15").Value = 11.9545335
$ws.Range("N15").Value = 23.909067
$ws.Range("O15").Value = 0.4359947946767024
$ws.Range("P15").Value = 0.3400888529957002
$ws.Range("Q15").Value = 4.96938600788175
$ws.Range("R15").Value = 19.877544031527
$ws.Range("S15").Value = 0.1067272690777811
$ws.Range("T15").Value = 0.06043128331517357

# Row 16
$ws.Range("G16").Value = 0.4156905
$ws.Range("H16").Value = 0.831381
$ws.Range("I16").Value = 0.2447902369039089
$ws.Range("J16").Value = 0.1776926317427335
$ws.Range("M16").Value = 0.5368526666666666
$ws.Range("N16").Value = 1.610558
$ws.Range("O16").Value = 0.01957959866647858
$ws.Range("P16").Value = 0.022909000292778
$ws.Range("Q16").Value = 0.223164553433
$ws.Range("R16").Value = 1.338987320598
$ws.Range("S16").Value = 0.004792894596050748
$ws.Range("T16").Value = 0.004070760552618775

